# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
# Both sheets share the same layout; update column F for rows 3, 4, 6, 7, 8.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    3 = 2499
    4 = 483
    6 = 6538
    7 = 371
    8 = 5
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
